# Update folder structure strings: insert "/Blocks" segment into the
# ParkDetails/InAroundPark/Components/ContentBlockComponent paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A52").Value = "ParkDetails/InAroundPark/Components/ContentBlockComponent/Blocks"
$ws.Range("A53").Value = "ParkDetails/InAroundPark/Components/ContentBlockComponent/Blocks/Title"
$ws.Range("A54").Value = "ParkDetails/InAroundPark/Components/ContentBlockComponent/Blocks/Text"
$ws.Range("A55").Value = "ParkDetails/InAroundPark/Components/ContentBlockComponent/Blocks/ImageContentid"
$ws.Range("A56").Value = "ParkDetails/InAroundPark/Components/ContentBlockComponent/Blocks/ImageAltText"

# Update the view state to match where the edit was made: scrolled down so
# row 43 is the top visible row, with A56 as the active/selected cell.
$ws.Activate()
$ws.Range("A56").Select()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
